$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.659.63'
$ws.Range("E2").Value = '  -3.19%  '
$ws.Range("D3").Value = '1.851.52'
$ws.Range("E3").Value = '  -3.88%  '
$ws.Range("E4").Value = '  -1.03%  '
$ws.Range("D5").Value = "'335.24"
$ws.Range("E5").Value = '  +2.75%  '
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D7").Value = "'0.4656"
$ws.Range("E7").Value = '  -3.39%  '
$ws.Range("D8").Value = "'0.3910"
$ws.Range("E8").Value = '  -3.69%  '
$ws.Range("D9").Value = "'46.22"
$ws.Range("E9").Value = '  -3.27%  '
$ws.Range("D10").Value = "'0.07905"
$ws.Range("E10").Value = '  -4.09%  '
$ws.Range("D11").Value = "'0.9851"
$ws.Range("E11").Value = '  -2.59%  '
$ws.Range("D12").Value = "'22.31"
$ws.Range("E12").Value = '  -6.59%  '
$ws.Range("D13").Value = '1.887.36'
$ws.Range("E13").Value = '  -1.64%  '
$ws.Range("D14").Value = "'5.855"
$ws.Range("E14").Value = '  -3.93%  '
$ws.Range("D15").Value = "'7.014"
$ws.Range("E15").Value = '  -3.27%  '
$ws.Range("D16").Value = "'0.06843"
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = "'87.79"
$ws.Range("E17").Value = '  -4.40%  '
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").Value = "'1.003"
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("E19").Value = '  -3.18%  '
$ws.Range("D20").Value = "'17.15"
$ws.Range("E20").Value = '  -2.66%  '
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D22").Value = '28.658.44'
$ws.Range("E22").Value = '  -3.17%  '
$ws.Range("D23").Value = "'5.405"
$ws.Range("E23").Value = '  -4.90%  '
$ws.Range("E24").Value = '  -5.87%  '
$ws.Range("D25").Value = "'2.136"
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("D26").Value = '2.071.81'
$ws.Range("E26").Value = '  -4.25%  '
$ws.Range("D27").Value = "'153.16"
$ws.Range("E27").Value = '  -2.00%  '
$ws.Range("D28").Value = "'19.52"
$ws.Range("D29").Value = "'6.055"
$ws.Range("E29").Value = '  -6.17%  '
$ws.Range("D30").Value = "'2.029"
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("D31").Value = "'117.78"
$ws.Range("D32").Value = "'0.9778"
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("D33").Value = "'0.09432"
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("D34").Value = "'5.380"
$ws.Range("E34").Value = '  -4.49%  '
$ws.Range("D35").Value = "'3.485"
$ws.Range("E35").Value = '  -2.69%  '
$ws.Range("D36").Value = "'1.353"
$ws.Range("D37").Value = "'0.06186"
$ws.Range("E37").Value = '  -3.35%  '
$ws.Range("D38").Value = "'0.02201"
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("D39").Value = "'1.166"
$ws.Range("E39").Value = '  -1.93%  '
$ws.Range("E40").Value = '  -0.93%  '
$ws.Range("D41").Value = "'0.5735"
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("D42").Value = "'7.634"
$ws.Range("E42").Value = '  -3.15%  '
$ws.Range("D43").Value = "'10.22"
$ws.Range("E43").Value = '  -4.93%  '
$ws.Range("D44").Value = "'0.1804"
$ws.Range("E44").Value = '  -2.60%  '
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("E46").Value = '  -3.23%  '
$ws.Range("D47").Value = "'0.5399"
$ws.Range("E47").Value = '  -2.99%  '
$ws.Range("D48").Value = "'11.78"
$ws.Range("E48").Value = '  -4.95%  '
$ws.Range("D49").Value = "'0.07149"
$ws.Range("E49").Value = '  -5.16%  '
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("D51").Value = "'114.95"
$ws.Range("E51").Value = '  -3.56%  '
